$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D7:L7")
$rng.NumberFormat = "[$-409]d\-mmm\-yy;@"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$ws.Range("D7").Value = 43373
$ws.Range("E7").Value = 43281
$ws.Range("F7").Value = 43190
$ws.Range("G7").Value = 43100
$ws.Range("H7").Value = 43008
$ws.Range("I7").Value = 42916
$ws.Range("J7").Value = 42825
$ws.Range("K7").Value = 42735
$ws.Range("L7").Value = 42643

$rng = $ws.Range("D8:L8")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D8").Value = 22100
$ws.Range("E8").Value = 20400
$ws.Range("F8").Value = 18700
$ws.Range("G8").Value = 18600
$ws.Range("H8").Value = 17400
$ws.Range("I8").Value = 15800
$ws.Range("J8").Value = 14600
$ws.Range("K8").Value = 13800
$ws.Range("L8").Value = 36900

$rng = $ws.Range("D9:L9")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "NA"
$ws.Range("H9").Value = "NA"
$ws.Range("I9").Value = "NA"
$ws.Range("J9").Value = "NA"
$ws.Range("K9").Value = "NA"
$ws.Range("L9").Value = "NA"

$rng = $ws.Range("D10:L10")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "NA"
$ws.Range("H10").Value = "NA"
$ws.Range("I10").Value = "NA"
$ws.Range("J10").Value = "NA"
$ws.Range("K10").Value = "NA"
$ws.Range("L10").Value = "NA"

$rng = $ws.Range("D11:L11")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152

$rng = $ws.Range("D12:L12")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"

$rng = $ws.Range("D13:L13")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0

$rng = $ws.Range("D14:L14")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0

$rng = $ws.Range("D15:L15")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0

$rng = $ws.Range("D16:L16")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152

$rng = $ws.Range("D17:L17")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D17").Value = 6800
$ws.Range("E17").Value = 5400
$ws.Range("F17").Value = 4500
$ws.Range("G17").Value = 4900
$ws.Range("H17").Value = 4500
$ws.Range("I17").Value = 3500
$ws.Range("J17").Value = 3400
$ws.Range("K17").Value = 3300
$ws.Range("L17").Value = 8500

$rng = $ws.Range("D18:L18")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D18").Value = 15300
$ws.Range("E18").Value = 15000
$ws.Range("F18").Value = 14200
$ws.Range("G18").Value = 13700
$ws.Range("H18").Value = 12900
$ws.Range("I18").Value = 12300
$ws.Range("J18").Value = 11200
$ws.Range("K18").Value = 10500
$ws.Range("L18").Value = 28400

$rng = $ws.Range("D19:L19")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152

$rng = $ws.Range("D20:L20")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D20").Value = -6700
$ws.Range("E20").Value = -6000
$ws.Range("F20").Value = -6100
$ws.Range("G20").Value = -8300
$ws.Range("H20").Value = -5100
$ws.Range("I20").Value = -4800
$ws.Range("J20").Value = -4800
$ws.Range("K20").Value = -5100
$ws.Range("L20").Value = -12500

$rng = $ws.Range("D21:L21")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D21").Value = 8900
$ws.Range("E21").Value = 9300
$ws.Range("F21").Value = 8300
$ws.Range("G21").Value = 5600
$ws.Range("H21").Value = 8000
$ws.Range("I21").Value = 7700
$ws.Range("J21").Value = 6700
$ws.Range("K21").Value = 5600
$ws.Range("L21").Value = 16300

$rng = $ws.Range("D22:L22")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0

$rng = $ws.Range("D23:L23")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D23").Value = 8600
$ws.Range("E23").Value = 9000
$ws.Range("F23").Value = 8000
$ws.Range("G23").Value = 5400
$ws.Range("H23").Value = 7700
$ws.Range("I23").Value = 7500
$ws.Range("J23").Value = 6500
$ws.Range("K23").Value = 5400
$ws.Range("L23").Value = 15900

$rng = $ws.Range("D24:L24")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D24").Value = 2200
$ws.Range("E24").Value = 2300
$ws.Range("F24").Value = 2100
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 3100
$ws.Range("I24").Value = 2700
$ws.Range("J24").Value = 2400
$ws.Range("K24").Value = 2000
$ws.Range("L24").Value = 6000

$rng = $ws.Range("D25:L25")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0

$rng = $ws.Range("D26:L26")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D26").Value = 6500
$ws.Range("E26").Value = 6700
$ws.Range("F26").Value = 6000
$ws.Range("G26").Value = 5300
$ws.Range("H26").Value = 4700
$ws.Range("I26").Value = 4800
$ws.Range("J26").Value = 4100
$ws.Range("K26").Value = 3400
$ws.Range("L26").Value = 9800

$rng = $ws.Range("D27:L27")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D27").Value = 6500
$ws.Range("E27").Value = 6700
$ws.Range("F27").Value = 6000
$ws.Range("G27").Value = 5300
$ws.Range("H27").Value = 4700
$ws.Range("I27").Value = 4800
$ws.Range("J27").Value = 4100
$ws.Range("K27").Value = 3400
$ws.Range("L27").Value = 9800

$rng = $ws.Range("D28:L28")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0

$rng = $ws.Range("D29:L29")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = -2000
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = "NA"
$ws.Range("L29").Value = "NA"

$rng = $ws.Range("D30:L30")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0

$rng = $ws.Range("D31:L31")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0

$rng = $ws.Range("D32:L32")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D32").Value = 6700
$ws.Range("E32").Value = 6000
$ws.Range("F32").Value = 6100
$ws.Range("G32").Value = 8300
$ws.Range("H32").Value = 5100
$ws.Range("I32").Value = 4800
$ws.Range("J32").Value = 4800
$ws.Range("K32").Value = 5100
$ws.Range("L32").Value = 12500

$rng = $ws.Range("D33:L33")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D33").Value = 6500
$ws.Range("E33").Value = 6700
$ws.Range("F33").Value = 6000
$ws.Range("G33").Value = 3300
$ws.Range("H33").Value = 4700
$ws.Range("I33").Value = 4800
$ws.Range("J33").Value = 4100
$ws.Range("K33").Value = 3400
$ws.Range("L33").Value = 9800

$rng = $ws.Range("D34:L34")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0

$rng = $ws.Range("D35:L35")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D35").Value = 6500
$ws.Range("E35").Value = 6700
$ws.Range("F35").Value = 6000
$ws.Range("G35").Value = 3300
$ws.Range("H35").Value = 4700
$ws.Range("I35").Value = 4800
$ws.Range("J35").Value = 4100
$ws.Range("K35").Value = 3400
$ws.Range("L35").Value = 9800

$rng = $ws.Range("D38:L38")
$rng.NumberFormat = "[$-409]d\-mmm\-yy;@"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$ws.Range("D38").Value = 43373
$ws.Range("E38").Value = 43281
$ws.Range("F38").Value = 43190
$ws.Range("G38").Value = 43100
$ws.Range("H38").Value = 43008
$ws.Range("I38").Value = 42916
$ws.Range("J38").Value = 42825
$ws.Range("K38").Value = 42735
$ws.Range("L38").Value = 42643

$rng = $ws.Range("D39:L39")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152

$rng = $ws.Range("D40:L40")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152

$rng = $ws.Range("D41:L41")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D41").Value = 21500
$ws.Range("E41").Value = 21900
$ws.Range("F41").Value = 20100
$ws.Range("G41").Value = 23700
$ws.Range("H41").Value = 24400
$ws.Range("I41").Value = "NA"
$ws.Range("J41").Value = "NA"
$ws.Range("K41").Value = "NA"
$ws.Range("L41").Value = "NA"

$rng = $ws.Range("D42:L42")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D42").Value = 11100
$ws.Range("E42").Value = 9100
$ws.Range("F42").Value = 9000
$ws.Range("G42").Value = 8200
$ws.Range("H42").Value = 7700
$ws.Range("I42").Value = "NA"
$ws.Range("J42").Value = "NA"
$ws.Range("K42").Value = "NA"
$ws.Range("L42").Value = "NA"

$rng = $ws.Range("D43:L43")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0

$rng = $ws.Range("D44:L44")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0

$rng = $ws.Range("D45:L45")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0

$rng = $ws.Range("D46:L46")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0

$rng = $ws.Range("D47:L47")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0

$rng = $ws.Range("D48:L48")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D48").Value = 11400
$ws.Range("E48").Value = 10500
$ws.Range("F48").Value = 10200
$ws.Range("G48").Value = 10100
$ws.Range("H48").Value = 10200
$ws.Range("I48").Value = "NA"
$ws.Range("J48").Value = "NA"
$ws.Range("K48").Value = "NA"
$ws.Range("L48").Value = "NA"

$rng = $ws.Range("D49:L49")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D49").Value = 3700
$ws.Range("E49").Value = 3800
$ws.Range("F49").Value = 3800
$ws.Range("G49").Value = 3900
$ws.Range("H49").Value = 3900
$ws.Range("I49").Value = "NA"
$ws.Range("J49").Value = "NA"
$ws.Range("K49").Value = "NA"
$ws.Range("L49").Value = "NA"

$rng = $ws.Range("D50:L50")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0

$rng = $ws.Range("D51:L51")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0

$rng = $ws.Range("D52:L52")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0

$rng = $ws.Range("D53:L53")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0

$rng = $ws.Range("D54:L54")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D54").Value = 1885800
$ws.Range("E54").Value = 1752900
$ws.Range("F54").Value = 1681600
$ws.Range("G54").Value = 1616600
$ws.Range("H54").Value = 1556700
$ws.Range("I54").Value = "NA"
$ws.Range("J54").Value = "NA"
$ws.Range("K54").Value = "NA"
$ws.Range("L54").Value = "NA"

$rng = $ws.Range("D55:L55")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152

$rng = $ws.Range("D56:L56")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152

$rng = $ws.Range("D57:L57")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D57").Value = 1200
$ws.Range("E57").Value = 1500
$ws.Range("F57").Value = 1100
$ws.Range("G57").Value = 1400
$ws.Range("H57").Value = 1000
$ws.Range("I57").Value = "NA"
$ws.Range("J57").Value = "NA"
$ws.Range("K57").Value = "NA"
$ws.Range("L57").Value = "NA"

$rng = $ws.Range("D58:L58")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0

$rng = $ws.Range("D59:L59")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0

$rng = $ws.Range("D60:L60")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0

$rng = $ws.Range("D61:L61")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D61").Value = 40100
$ws.Range("E61").Value = 40600
$ws.Range("F61").Value = 41100
$ws.Range("G61").Value = 41500
$ws.Range("H61").Value = 42000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0

$rng = $ws.Range("D62:L62")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0

$rng = $ws.Range("D63:L63")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0

$rng = $ws.Range("D64:L64")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0

$rng = $ws.Range("D65:L65")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0

$rng = $ws.Range("D66:L66")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D66").Value = 1674900
$ws.Range("E66").Value = 1547000
$ws.Range("F66").Value = 1482600
$ws.Range("G66").Value = 1479500
$ws.Range("H66").Value = 1423000
$ws.Range("I66").Value = "NA"
$ws.Range("J66").Value = "NA"
$ws.Range("K66").Value = "NA"
$ws.Range("L66").Value = "NA"

$rng = $ws.Range("D67:L67")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152

$rng = $ws.Range("D68:L68")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0

$rng = $ws.Range("D69:L69")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0

$rng = $ws.Range("D70:L70")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0

$rng = $ws.Range("D71:L71")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0

$rng = $ws.Range("D72:L72")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D72").Value = 88500
$ws.Range("E72").Value = 82000
$ws.Range("F72").Value = 75300
$ws.Range("G72").Value = 69500
$ws.Range("H72").Value = 66200
$ws.Range("I72").Value = "NA"
$ws.Range("J72").Value = "NA"
$ws.Range("K72").Value = "NA"
$ws.Range("L72").Value = "NA"

$rng = $ws.Range("D73:L73")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0

$rng = $ws.Range("D74:L74")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0

$rng = $ws.Range("D75:L75")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0

$rng = $ws.Range("D76:L76")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D76").Value = 210900
$ws.Range("E76").Value = 205900
$ws.Range("F76").Value = 199000
$ws.Range("G76").Value = 137200
$ws.Range("H76").Value = 133700
$ws.Range("I76").Value = "NA"
$ws.Range("J76").Value = "NA"
$ws.Range("K76").Value = "NA"
$ws.Range("L76").Value = "NA"

$rng = $ws.Range("D77:L77")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0

$rng = $ws.Range("D80:L80")
$rng.NumberFormat = "[$-409]d\-mmm\-yy;@"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$ws.Range("D80").Value = 43373
$ws.Range("E80").Value = 43281
$ws.Range("F80").Value = 43190
$ws.Range("G80").Value = 43100
$ws.Range("H80").Value = 43008
$ws.Range("I80").Value = 42916
$ws.Range("J80").Value = 42825
$ws.Range("K80").Value = 42735
$ws.Range("L80").Value = 42643

$rng = $ws.Range("D81:L81")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D81").Value = 6500
$ws.Range("E81").Value = 6700
$ws.Range("F81").Value = 6000
$ws.Range("G81").Value = 3300
$ws.Range("H81").Value = 4700
$ws.Range("I81").Value = 4800
$ws.Range("J81").Value = 4100
$ws.Range("K81").Value = 3400
$ws.Range("L81").Value = 9800

$rng = $ws.Range("D82:L82")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152

$rng = $ws.Range("D83:L83")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D83").Value = 200
$ws.Range("E83").Value = 200
$ws.Range("F83").Value = 200
$ws.Range("G83").Value = 200
$ws.Range("H83").Value = 200
$ws.Range("I83").Value = 200
$ws.Range("J83").Value = 200
$ws.Range("K83").Value = 200
$ws.Range("L83").Value = 500

$rng = $ws.Range("D84:L84")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0

$rng = $ws.Range("D85:L85")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0

$rng = $ws.Range("D86:L86")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0

$rng = $ws.Range("D87:L87")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0

$rng = $ws.Range("D88:L88")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0

$rng = $ws.Range("D89:L89")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D89").Value = 9800
$ws.Range("E89").Value = 7900
$ws.Range("F89").Value = 4500
$ws.Range("G89").Value = 5900
$ws.Range("H89").Value = 6600
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 8500
$ws.Range("K89").Value = 2400
$ws.Range("L89").Value = 12600

$rng = $ws.Range("D90:L90")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152

$rng = $ws.Range("D91:L91")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D91").Value = -1100
$ws.Range("E91").Value = -500
$ws.Range("F91").Value = -200
$ws.Range("G91").Value = -100
$ws.Range("H91").Value = -100
$ws.Range("I91").Value = -300
$ws.Range("J91").Value = -700
$ws.Range("K91").Value = -1000
$ws.Range("L91").Value = -1200

$rng = $ws.Range("D92:L92")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0

$rng = $ws.Range("D93:L93")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0

$rng = $ws.Range("D94:L94")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D94").Value = -137100
$ws.Range("E94").Value = -69200
$ws.Range("F94").Value = -71200
$ws.Range("G94").Value = -62800
$ws.Range("H94").Value = -106800
$ws.Range("I94").Value = -85400
$ws.Range("J94").Value = -95300
$ws.Range("K94").Value = -79300
$ws.Range("L94").Value = -186800

$rng = $ws.Range("D95:L95")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152

$rng = $ws.Range("D96:L96")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0

$rng = $ws.Range("D97:L97")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0

$rng = $ws.Range("D98:L98")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0

$rng = $ws.Range("D99:L99")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0

$rng = $ws.Range("D100:L100")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D100").Value = 126900
$ws.Range("E100").Value = 63200
$ws.Range("F100").Value = 63000
$ws.Range("G100").Value = 56200
$ws.Range("H100").Value = 102800
$ws.Range("I100").Value = 86700
$ws.Range("J100").Value = 86800
$ws.Range("K100").Value = 79400
$ws.Range("L100").Value = 174700

$rng = $ws.Range("D101:L101")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0

$rng = $ws.Range("D102:L102")
$rng.NumberFormat = "#,##0"
$rng.Font.Name = "Verdana"
$rng.Font.Size = 12
$rng.Font.Bold = $false
$rng.HorizontalAlignment = -4152
$ws.Range("D102").Value = -400
$ws.Range("E102").Value = 1800
$ws.Range("F102").Value = -3600
$ws.Range("G102").Value = -700
$ws.Range("H102").Value = 2600
$ws.Range("I102").Value = 5300
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2500
$ws.Range("L102").Value = 400
